$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; this shifts the existing "District" column
# (old F) one position to the right, into the new column G.
$ws.Columns.Item(6).Insert()

# Set the new header for column F ("Address")
$ws.Range("F2").Value2 = "Address"

# Populate the new Address values in column F for each data row
$ws.Range("F3").Value2 = "G R H S KadapalakerePavagada"
$ws.Range("F4").Value2 = "G H SBrahmasandraSira"
$ws.Range("F5").Value2 = "G J CKoratagere"
$ws.Range("F6").Value2 = "M G M Girls High School"
$ws.Range("F7").Value2 = "Sri Gajanana Girls High School Sira"
$ws.Range("F8").Value2 = "G H S RajavanthiPavagad"
$ws.Range("F9").Value2 = "T R High SchoolMangalavadaPavagada"
$ws.Range("F10").Value2 = "G H S AkkirampuraKoratagere"
$ws.Range("F11").Value2 = "S L N High SchoolMedigeshi"
$ws.Range("F12").Value2 = "R P H S MugadalabettaPavagada"
$ws.Range("F13").Value2 = "G J CChikkanahalliSira"
$ws.Range("F14").Value2 = "G H S ArasapuraKoratagere"
$ws.Range("F15").Value2 = "G H S MAddakkanahalliSira"
$ws.Range("F16").Value2 = "G J C Koratagere"
$ws.Range("F17").Value2 = "S R R H S Neralekere"
$ws.Range("F18").Value2 = "G H S BaraguruSira"
$ws.Range("F19").Value2 = "S B R H S ShidlekeneSira"
$ws.Range("F20").Value2 = "S G H S Kodigenahally"
$ws.Range("F21").Value2 = "R R H S Hosakere"
$ws.Range("F22").Value2 = "G J C LakkanahalliSira"
$ws.Range("F23").Value2 = "G H S Krishnapura"
$ws.Range("F25").Value2 = "G G H SKoratagere"
$ws.Range("F26").Value2 = "G H S Holavana halliKoratagere"
$ws.Range("F27").Value2 = "G H S VaddagereKoratagere"
$ws.Range("F28").Value2 = "G J C GuligenahalliSira"
$ws.Range("F29").Value2 = "G J C (High School Section) H L DurgaKunigal"
$ws.Range("F30").Value2 = "Govt. High SchoolMavinakereTuruvekere"
$ws.Range("F31").Value2 = "G J C (High School section) B H Road"
$ws.Range("F32").Value2 = "Govt. High School IppadiKunigal"
$ws.Range("F33").Value2 = "K M H P S C N Hally"
$ws.Range("F34").Value2 = "G H P S KichchavadiKunigal"
$ws.Range("F35").Value2 = "G M H P S SantemavathurKunigal"
$ws.Range("F36").Value2 = "G H P S H ThammadihallyC N Hally"
$ws.Range("F37").Value2 = "G B H S ChelurGubbi"
$ws.Range("F38").Value2 = "R R H S Ankanahalli MuttKunigal"
$ws.Range("F39").Value2 = "G H S YalanaduC N Hally"
$ws.Range("F40").Value2 = "S P S J C (H S) RangapuraTiptur"
$ws.Range("F41").Value2 = "G H P S MarendupalyaC N Halli"
$ws.Range("F42").Value2 = "G H P S ChikkarampuraC N Halli"
$ws.Range("F43").Value2 = "S G H S Vijaya nagar"
$ws.Range("F44").Value2 = "S L B S H S SarthavallyTiptur"
$ws.Range("F45").Value2 = "G M P S YalanaduC N Hally"
$ws.Range("F46").Value2 = "P T G R M H STiptur"
$ws.Range("F47").Value2 = "N S M G H STiptur"
$ws.Range("F48").Value2 = "G H P S HonnebagiC N Hally"
$ws.Range("F49").Value2 = "G H S GanganagattaTiptur"
$ws.Range("F50").Value2 = "G H S NellikereTiptur"
$ws.Range("F51").Value2 = "G H S DabbegattaTuruvekere"
